$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.446.44'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.80%  '
$ws.Range('E2').ClearFormats()

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.532.74'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.77%  '
$ws.Range('E3').ClearFormats()

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('E4').ClearFormats()

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '611.84'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -5.05%  '
$ws.Range('E5').ClearFormats()

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.86'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -3.41%  '
$ws.Range('E6').ClearFormats()

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.527.19'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -3.88%  '
$ws.Range('E7').ClearFormats()

# Row 8
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('E8').ClearFormats()

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.487'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.77%  '
$ws.Range('E9').ClearFormats()

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.141'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.37%  '
$ws.Range('E10').ClearFormats()

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.91'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.84%  '
$ws.Range('E11').ClearFormats()

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.430'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.40%  '
$ws.Range('E12').ClearFormats()

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000222'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -4.09%  '
$ws.Range('E13').ClearFormats()

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.133.31'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.69%  '
$ws.Range('E14').ClearFormats()

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '32.01'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.67%  '
$ws.Range('E15').ClearFormats()

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.551.29'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.19%  '
$ws.Range('E16').ClearFormats()

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.532.96'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.69%  '
$ws.Range('E17').ClearFormats()

# Row 18
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.84%  '
$ws.Range('E18').ClearFormats()

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.38'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.42%  '
$ws.Range('E19').ClearFormats()

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.49'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.15%  '
$ws.Range('E20').ClearFormats()

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '454.61'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.30%  '
$ws.Range('E21').ClearFormats()

# Row 22
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -4.57%  '
$ws.Range('E22').ClearFormats()

# Row 23
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.29%  '
$ws.Range('E23').ClearFormats()

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.68'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.96%  '
$ws.Range('E24').ClearFormats()

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.676.93'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.71%  '
$ws.Range('E25').ClearFormats()

# Row 26
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E26').ClearFormats()

# Row 27
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.97%  '
$ws.Range('E27').ClearFormats()

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.46'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.52%  '
$ws.Range('E28').ClearFormats()

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.36'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -6.63%  '
$ws.Range('E29').ClearFormats()

# Row 30
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.16%  '
$ws.Range('E30').ClearFormats()

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.68'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.92%  '
$ws.Range('E31').ClearFormats()

# Row 32
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('E32').ClearFormats()

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.04'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.25%  '
$ws.Range('E33').ClearFormats()

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.91'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.27%  '
$ws.Range('E34').ClearFormats()

# Row 35
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.159'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.46%  '
$ws.Range('E35').ClearFormats()

# Row 36
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.21'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.91%  '
$ws.Range('E36').ClearFormats()

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.532.35'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.54%  '
$ws.Range('E37').ClearFormats()

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.99'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.97%  '
$ws.Range('E38').ClearFormats()

# Row 39
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('E39').ClearFormats()

# Row 40
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('E40').ClearFormats()

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '175.92'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.11%  '
$ws.Range('E41').ClearFormats()

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.62'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -4.65%  '
$ws.Range('E42').ClearFormats()

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0880'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.99%  '
$ws.Range('E43').ClearFormats()

# Row 44
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.86%  '
$ws.Range('E44').ClearFormats()

# Row 45
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.66%  '
$ws.Range('E45').ClearFormats()

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '29.41'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +9.83%  '
$ws.Range('E46').ClearFormats()

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '45.75'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.92%  '
$ws.Range('E47').ClearFormats()

# Row 48
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.61%  '
$ws.Range('E48').ClearFormats()

# Row 49
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.64%  '
$ws.Range('E49').ClearFormats()

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.67'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.02%  '
$ws.Range('E50').ClearFormats()

# Row 51
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.17%  '
$ws.Range('E51').ClearFormats()
